$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ------------------------------------------------------------------
# 1. Append four new rows describing the new admin delete/feature
#    requests. Rows.Add() clones the formatting (shading, language,
#    justification) of the table's current last row, which is exactly
#    the shading/formatting these new rows need.
# ------------------------------------------------------------------

# --- Row: delete / Admin/deleteJobItem/:jobItemId/:id/:token ---
$r1 = $t.Rows.Add()
$r1.Cells.Item(1).Range.Text = "delete"
$r1.Cells.Item(2).Range.Text = "Admin/deleteJobItem/:jobItemId/:id/:token"
$r1.Cells.Item(5).Range.Text = "Delete the disered jobitem"

# --- Row: delete / Admin/deleteCompany/:companyId/:id/:token ---
$r2 = $t.Rows.Add()
$r2.Cells.Item(1).Range.Text = "delete"
$r2.Cells.Item(2).Range.Text = "Admin/deleteCompany/:companyId/:id/:token"
$r2.Cells.Item(5).Range.Text = "Delete the desired company user"

# --- Row: put / Admin/featureCompany/:companyId/:id/:token ---
$r3 = $t.Rows.Add()
$r3.Cells.Item(1).Range.Text = "put"
$r3.Cells.Item(2).Range.Text = "Admin/featureCompany/:companyId/:id/:token"
$r3.Cells.Item(3).Range.Text = "feature (true/false)"
$r3.Cells.Item(5).Range.Text = "Feature the desired company"

# --- Row: put / Admin/featureJobitem/:jobItemId/:id/:token ---
$r4 = $t.Rows.Add()
$r4.Cells.Item(1).Range.Text = "put"
$r4.Cells.Item(2).Range.Text = "Admin/featureJobitem/:jobItemId/:id/:token"
$r4.Cells.Item(3).Range.Text = "feature (true/false)"
$r4.Cells.Item(5).Range.Text = "Feature the desired jobitem"

# ------------------------------------------------------------------
# 2. The extra rows/content made column 2 ("path") need a little more
#    room, so the whole grid was rebalanced (same overall table width)
#    -- shrinking columns 3/4/5 slightly to grow column 2.
# ------------------------------------------------------------------
$t.Columns.Item(2).Width = 224.65   # 4493 dxa
$t.Columns.Item(3).Width = 158.1    # 3162 dxa
$t.Columns.Item(4).Width = 91.6     # 1832 dxa
$t.Columns.Item(5).Width = 252.8    # 5056 dxa
